# Add a new "Turkey" sheet, cloned from the existing "Spain" sheet (same
# layout/styles), populate it with the Turkey-specific market values, and
# tidy up the selection/view state on both the source and the new sheet -
# mirroring what Excel does when you right-click a tab -> "Move or Copy..."
# -> Create a copy, then edit the two data cells.

$wb = $excel.ActiveWorkbook

$spain = $wb.Worksheets.Item("Spain")

# Duplicate "Spain" and drop the copy right after it (i.e. at the end of
# the tab strip), exactly like Excel's own "(Move or) Copy" does.
$spain.Copy($null, $spain)
$turkey = $wb.Worksheets.Item($wb.Worksheets.Count)
$turkey.Name = "Turkey"

# Turkey-specific data (the two cells that differ from the Spain template).
$turkey.Range("B2").Value = "Turkey Market"
$turkey.Range("B4").Value = "NGC-3191/T3313"

# The longer "NGC-3191/T3313" text plus the shorter "Turkey Market" causes
# column D to shrink back toward the sheet's default width, with rows 3-5
# wrapping onto two lines - reproduce that layout explicitly.
$turkey.Columns.Item(4).ColumnWidth = 7.6
$turkey.Rows.Item(3).RowHeight = 28.8
$turkey.Rows.Item(4).RowHeight = 28.8
$turkey.Rows.Item(5).RowHeight = 28.8

# Clear the old selection/active-tab state left over on "Spain" now that
# focus has moved to the new sheet, and select the new sheet's last-used
# cell so it becomes the active tab.
$spain.Select()
$spain.Range("A1:D12").Select()

$turkey.Select()
$turkey.Range("I16").Select()
